$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("currentshop")

$data = @(
    @("Taco Hat", 60, 31, 1),
    @("Gold Spatula", 100, 8, 1),
    @("Gold Knife", 100, 9, 1),
    @("Sombrero", 120, 32, 2),
    @("Extra Burner", 150, 28, 3),
    @("Extra Burner 2", 200, 29, 3),
    @("Small Plant", 13, 16, 4),
    @("Hard Taco Sign", 25, 4, 4),
    @("Soft Taco Sign", 25, 5, 4),
    @("Pita Sign", 25, 6, 4),
    @("Chef Hat", 180, 33, 5),
    @("Romano Poster", 50, 12, 6),
    @("Viking Helmet", 300, 34, 6),
    @("Beans Poster", 8, 3, 7),
    @("Cactus", 10, 2, 7),
    @("Mystery Poster", 20, 10, 7),
    @("Burgeria Poster", 20, 11, 7),
    @("Desert Plant", 25, 15, 7),
    @("Leafy Plant", 40, 17, 7),
    @("Onion Ring Poster", 50, 13, 7),
    @("Pepper Lights", 60, 7, 7),
    @("Coffee Stand", 65, 23, 7),
    @("Newspaper Stand", 80, 19, 7),
    @("Hypno Clock", 80, 30, 7),
    @("Gumball Machine", 120, 21, 7),
    @("Ceiling Fan", 150, 14, 7),
    @("TV", 150, 18, 7),
    @("Arcade Cabinet", 400, 20, 7),
    @("Jukebox", 500, 22, 7),
    @("Royal Crown", 1000, 35, 7),
    @("Doorbell", 30, 1, 8),
    @("Beef Alarm", 90, 24, 8),
    @("Chicken Alarm", 90, 25, 8),
    @("Pork Alarm", 90, 26, 8),
    @("Steak Alarm", 90, 27, 8)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r++
}
